$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44432
$ws.Range("M2").Value = 30
$ws.Range("N2").Value = 1300
$ws.Range("O2").Value = 1300
$ws.Range("P2").Value = 1300
$ws.Range("S2").Value = 1300

# Row 4
$ws.Range("D4").Value = 44431
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 1300
$ws.Range("O4").Value = 1300
$ws.Range("P4").Value = 1300
$ws.Range("S4").Value = 1300

# Row 5
$ws.Range("D5").Value = 44435
$ws.Range("M5").Value = 130

# Row 6
$ws.Range("D6").Value = 44417
$ws.Range("M6").Value = 80

# Row 7
$ws.Range("D7").Value = 44405
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 1200
$ws.Range("O7").Value = 1200
$ws.Range("P7").Value = 1200
$ws.Range("S7").Value = 1200

# Row 8
$ws.Range("D8").Value = 44418
$ws.Range("M8").Value = 40

# Row 9
$ws.Range("D9").Value = 44343
$ws.Range("M9").Value = 60

# Row 11
$ws.Range("D11").Value = 44357
$ws.Range("M11").Value = 35
$ws.Range("N11").Value = 1000
$ws.Range("O11").Value = 1000
$ws.Range("P11").Value = 1000
$ws.Range("S11").Value = 1000
